# Update the "Read Connections" example workbook.
# The original sheet had a stray header row (Cl1 / Node1) in row 1 that
# doesn't belong to the actual connection matrix, which starts on what
# was row 2. Remove that first row so the matrix shifts up by one row
# and select the new (now empty) top row, mirroring how it was left in
# Excel before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire first row (this is how the sheet was left selected
# in the authored file right before the row got removed).
$ws.Range("A1:XFD1").Select() | Out-Null

# Delete row 1 entirely; everything below shifts up by one row.
$ws.Rows("1:1").Delete() | Out-Null
